# update get_tweets without tweepy
#
# For every data row (2-33) on Sheet1:
#   - column E (last_id_tweet) goes from "1" to "0"
#   - column F (date_tweet) is rewritten to a sequential run of synthetic
#     timestamps starting at 2010-05-16 21:50:00 (one second apart per row),
#     replacing the old scraped tweet timestamps.
# Columns A (date_last_update), B (timestamp_last_update), C (id_user) and
# D (name_user) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 33; $r++) {
    $offset = $r - 2
    $ws.Cells.Item($r, 5).Value = "0"
    $ws.Cells.Item($r, 6).Value = ("2010-05-16 21:50:{0:D2}" -f $offset)
}

# Refresh the active selection/view to match the post-edit cursor position.
$ws.Range("F10").Select()
